# Apply updated crypto price/volume figures scraped on 2023-11-03.
# D-column ("Price") cells are forced to Text format before the write so
# Excel stores the decimal-looking strings verbatim instead of silently
# re-interpreting them as numbers (which would drop meaningful trailing
# zeros / separators, e.g. "0.960" -> 0.96). Style is reset to Normal right
# after so the cell keeps the workbook default formatting.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "34.629.53"
$ws.Range("E2").Value = "  -2.28%  "
Set-TextValue "D3" "1.806.76"
$ws.Range("E3").Value = "  -1.71%  "
$ws.Range("E4").Value = "  +0.34%  "
Set-TextValue "D5" "230.23"
$ws.Range("E5").Value = "  -0.50%  "
Set-TextValue "D6" "0.612"
$ws.Range("E6").Value = "  +0.35%  "
$ws.Range("E7").Value = "  +0.43%  "
Set-TextValue "D8" "39.33"
$ws.Range("E8").Value = "  -9.79%  "
$ws.Range("E9").Value = "  +4.64%  "
Set-TextValue "D10" "0.0681"
$ws.Range("E10").Value = "  -3.14%  "
$ws.Range("E11").Value = "  -1.78%  "
Set-TextValue "D12" "2.069.46"
$ws.Range("E12").Value = "  -1.69%  "
Set-TextValue "D13" "11.22"
$ws.Range("E13").Value = "  -0.24%  "
$ws.Range("E14").Value = "  -1.39%  "
Set-TextValue "D15" "1.799.30"
$ws.Range("E15").Value = "  -2.14%  "
Set-TextValue "D16" "4.59"
$ws.Range("E16").Value = "  -2.51%  "
Set-TextValue "D17" "34.645.31"
$ws.Range("E17").Value = "  -2.11%  "
Set-TextValue "D18" "69.41"
$ws.Range("E18").Value = "  -0.88%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
Set-TextValue "D19" "241.92"
$ws.Range("E19").Value = "  -0.97%  "
$ws.Range("B20").Value = "ShibaInu"
$ws.Range("C20").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
Set-TextValue "D20" "0.0₃0783"
$ws.Range("E20").Value = "  -2.11%  "
Set-TextValue "D21" "11.89"
$ws.Range("E21").Value = "  -1.77%  "
$ws.Range("E22").Value = "  -0.58%  "
$ws.Range("E23").Value = "  +0.45%  "
Set-TextValue "D24" "2.24"
$ws.Range("E24").Value = "  +1.64%  "
Set-TextValue "D25" "171.85"
$ws.Range("E25").Value = "  +0.67%  "
Set-TextValue "D26" "7.77"
$ws.Range("E26").Value = "  -1.75%  "
Set-TextValue "D27" "17.22"
$ws.Range("E27").Value = "  -3.04%  "
$ws.Range("E28").Value = "  +0.21%  "
Set-TextValue "D29" "1.48"
$ws.Range("E29").Value = "  -4.71%  "
$ws.Range("E30").Value = "  +0.30%  "
Set-TextValue "D31" "4.06"
$ws.Range("E31").Value = "  +3.20%  "
$ws.Range("E32").Value = "  -1.44%  "
$ws.Range("E33").Value = "  -4.01%  "
$ws.Range("E34").Value = "  +13.94%  "
$ws.Range("E35").Value = "  -3.27%  "
Set-TextValue "D36" "0.698"
$ws.Range("E36").Value = "  +1.52%  "
Set-TextValue "D37" "91.55"
$ws.Range("E37").Value = "  -4.50%  "
$ws.Range("E38").Value = "  +4.60%  "
Set-TextValue "D39" "1.325.43"
$ws.Range("E39").Value = "  -1.68%  "
$ws.Range("E40").Value = "  -1.94%  "
$ws.Range("E41").Value = "  +0.31%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
Set-TextValue "D42" "0.960"
$ws.Range("E42").Value = "  -5.08%  "
$ws.Range("B43").Value = "InjectiveProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
Set-TextValue "D43" "14.34"
$ws.Range("E43").Value = "  -8.30%  "
Set-TextValue "D44" "2.73"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("E45").Value = "  -9.77%  "
$ws.Range("E46").Value = "  -0.70%  "
Set-TextValue "D47" "0.0513"
$ws.Range("E47").Value = "  -1.27%  "
Set-TextValue "D48" "1.993.82"
$ws.Range("E48").Value = "  -0.67%  "
$ws.Range("E49").Value = "  +0.44%  "
Set-TextValue "D50" "0.0663"
$ws.Range("E50").Value = "  +5.30%  "
Set-TextValue "D51" "98.01"
$ws.Range("E51").Value = "  -5.02%  "
